$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1300.2667
$ws.Range("J17").Value = 1300.2667
$ws.Range("L17").Value = 3900.800099999999
$ws.Range("N17").Value = -4236.800099999999

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H113").Value = 27790134
$ws.Range("J113").Value = 33349912
$ws.Range("L113").Value = 33349912
$ws.Range("N113").Value = -33356420

$ws.Range("H116").Value = 15635809
$ws.Range("I116").Value = 22735996
$ws.Range("J116").Value = 15397.8
$ws.Range("K116").Value = 22735996
$ws.Range("L116").Value = 15397.8
$ws.Range("M116").Value = -22732554
$ws.Range("N116").Value = -22281.8

$ws.Range("H125").Value = 62501452
$ws.Range("I125").Value = 90910220
$ws.Range("J125").Value = 2155.6
$ws.Range("K125").Value = 818191980
$ws.Range("L125").Value = 19400.4
$ws.Range("M125").Value = -818189520
$ws.Range("N125").Value = -24320.4

$ws.Range("H138").Value = 1591121.8
$ws.Range("I138").Value = 1804.1923
$ws.Range("K138").Value = 5412.5769
$ws.Range("M138").Value = -272.5769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4769119
$ws.Range("I32").Value = 5133269
$ws.Range("K32").Value = 5133269
$ws.Range("M32").Value = -5132982

$ws.Range("H63").Value = 1977.5
$ws.Range("I63").Value = 1983.3334
$ws.Range("J63").Value = 1968.75
$ws.Range("K63").Value = 1983.3334
$ws.Range("L63").Value = 1968.75
$ws.Range("M63").Value = -1297.3334
$ws.Range("N63").Value = -3340.75

$ws.Range("H66").Value = 1977.5
$ws.Range("I66").Value = 1983.3334
$ws.Range("J66").Value = 1968.75
$ws.Range("K66").Value = 9916.666999999999
$ws.Range("L66").Value = 9843.75
$ws.Range("M66").Value = -6484.666999999999
$ws.Range("N66").Value = -16707.75

$ws.Range("H74").Value = 34958.418
$ws.Range("I74").Value = 47040.816
$ws.Range("J74").Value = 5423.6665
$ws.Range("K74").Value = 47040.816
$ws.Range("L74").Value = 5423.6665
$ws.Range("M74").Value = -46166.816
$ws.Range("N74").Value = -7171.6665

$ws.Range("H77").Value = 34958.418
$ws.Range("I77").Value = 47040.816
$ws.Range("J77").Value = 5423.6665
$ws.Range("K77").Value = 235204.08
$ws.Range("L77").Value = 27118.3325
$ws.Range("M77").Value = -230836.08
$ws.Range("N77").Value = -35854.3325

$ws.Range("H102").Value = 3556.4285
$ws.Range("J102").Value = 4377.6665
$ws.Range("L102").Value = 4377.6665
$ws.Range("N102").Value = -7621.6665

$ws.Range("H122").Value = 4078.2083
$ws.Range("I122").Value = 1614.4546
$ws.Range("J122").Value = 6162.923
$ws.Range("K122").Value = 4843.3638
$ws.Range("L122").Value = 18488.769
$ws.Range("M122").Value = -2393.3638
$ws.Range("N122").Value = -23388.769

$ws.Range("H132").Value = 5457.8184
$ws.Range("I132").Value = 3968.9167
$ws.Range("K132").Value = 11906.7501
$ws.Range("M132").Value = -9376.750100000001

$ws.Range("H133").Value = 87654.336
$ws.Range("J133").Value = 87654.336
$ws.Range("L133").Value = 87654.336
$ws.Range("N133").Value = -92714.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 66971
$ws.Range("J27").Value = 66971
$ws.Range("L27").Value = 66971
$ws.Range("N27").Value = -67355

$ws.Range("H33").Value = 750
$ws.Range("I33").Value = 500
$ws.Range("K33").Value = 500
$ws.Range("M33").Value = -164

$ws.Range("H107").Value = 43274132
$ws.Range("I107").Value = 46880070
$ws.Range("J107").Value = 2845
$ws.Range("K107").Value = 46880070
$ws.Range("L107").Value = 2845
$ws.Range("M107").Value = -46878150
$ws.Range("N107").Value = -6685

$ws.Range("H130").Value = 73350.11
$ws.Range("J130").Value = 73350.11
$ws.Range("L130").Value = 73350.11
$ws.Range("N130").Value = -83390.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 28417508
$ws.Range("I86").Value = 34729290
$ws.Range("K86").Value = 34729290
$ws.Range("M86").Value = -34728167

$ws.Range("H89").Value = 28417508
$ws.Range("I89").Value = 34729290
$ws.Range("K89").Value = 173646450
$ws.Range("M89").Value = -173640834

$ws.Range("H99").Value = 6853.25
$ws.Range("I99").Value = 4656
$ws.Range("J99").Value = 7585.6665
$ws.Range("K99").Value = 4656
$ws.Range("L99").Value = 7585.6665
$ws.Range("M99").Value = -3158
$ws.Range("N99").Value = -10581.6665

$ws.Range("H126").Value = 6853.25
$ws.Range("I126").Value = 4656
$ws.Range("J126").Value = 7585.6665
$ws.Range("K126").Value = 13968
$ws.Range("L126").Value = 22756.9995
$ws.Range("M126").Value = -11498
$ws.Range("N126").Value = -27696.9995

$ws.Range("H132").Value = 6002.3
$ws.Range("J132").Value = 8506
$ws.Range("L132").Value = 25518
$ws.Range("N132").Value = -30578

$ws.Range("H141").Value = 434249.75
$ws.Range("J141").Value = 434249.75
$ws.Range("L141").Value = 434249.75
$ws.Range("N141").Value = -444609.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 198.6875
$ws.Range("I33").Value = 72.85714
$ws.Range("K33").Value = 437.14284
$ws.Range("M33").Value = -154.14284

$ws.Range("H131").Value = 3274.4062
$ws.Range("I131").Value = 2599.6667
$ws.Range("K131").Value = 7799.000100000001
$ws.Range("M131").Value = -2759.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 61761.94
$ws.Range("I57").Value = 3500
$ws.Range("J57").Value = 69530.2
$ws.Range("K57").Value = 3500
$ws.Range("L57").Value = 69530.2
$ws.Range("M57").Value = -2680
$ws.Range("N57").Value = -71170.2

$ws.Range("H86").Value = 50003
$ws.Range("J86").Value = 50003
$ws.Range("L86").Value = 50003
$ws.Range("N86").Value = -52375

$ws.Range("H89").Value = 50003
$ws.Range("J89").Value = 50003
$ws.Range("L89").Value = 150009
$ws.Range("N89").Value = -161865

$ws.Range("H132").Value = 3641.9714
$ws.Range("I132").Value = 1638.1852
$ws.Range("K132").Value = 4914.5556
$ws.Range("M132").Value = -2384.5556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7483.75
$ws.Range("J40").Value = 7483.75
$ws.Range("L40").Value = 7483.75
$ws.Range("N40").Value = -7755.75

$ws.Range("H56").Value = 39051
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680

$ws.Range("H109").Value = 58000
$ws.Range("J109").Value = 58000
$ws.Range("L109").Value = 58000
$ws.Range("N109").Value = -60774

$ws.Range("H122").Value = 3534.1428
$ws.Range("I122").Value = 2097.6667
$ws.Range("K122").Value = 6293.000100000001
$ws.Range("M122").Value = -3843.000100000001

$ws.Range("H132").Value = 6832.386
$ws.Range("I132").Value = 3539.0625
$ws.Range("K132").Value = 10617.1875
$ws.Range("M132").Value = -8087.1875

$ws.Range("H136").Value = 12257.447
$ws.Range("I136").Value = 3040
$ws.Range("J136").Value = 19085.186
$ws.Range("K136").Value = 9120
$ws.Range("L136").Value = 57255.558
$ws.Range("M136").Value = -6570
$ws.Range("N136").Value = -62355.558

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H107").Value = 10101818
$ws.Range("I107").Value = 382.26315
$ws.Range("J107").Value = 23810910
$ws.Range("K107").Value = 1146.78945
$ws.Range("L107").Value = 71432730
$ws.Range("M107").Value = 773.21055
$ws.Range("N107").Value = -71436570

$ws.Range("H122").Value = 2603.2195
$ws.Range("I122").Value = 2035.4138
$ws.Range("K122").Value = 6106.2414
$ws.Range("M122").Value = -3656.2414

$ws.Range("H125").Value = 49928
$ws.Range("J125").Value = 49928
$ws.Range("L125").Value = 49928
$ws.Range("N125").Value = -59768

$ws.Range("H126").Value = 1194.7693
$ws.Range("I126").Value = 1193
$ws.Range("K126").Value = 3579
$ws.Range("M126").Value = -1109

$ws.Range("H132").Value = 2800.0417
$ws.Range("I132").Value = 1199.8125
$ws.Range("K132").Value = 3599.4375
$ws.Range("M132").Value = -1069.4375

$ws.Range("H136").Value = 4490.75
$ws.Range("I136").Value = 1672.4
$ws.Range("K136").Value = 5017.200000000001
$ws.Range("M136").Value = -2467.200000000001

